$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows (dates in column A, hours worked in column B)
# Row 76: 13-Aug, 1 hour
$ws.Cells.Item(76, 1).Value = 45882
$ws.Cells.Item(76, 2).Value = 1

# Row 77: 22-Aug, 0.5 hour
$ws.Cells.Item(77, 1).Value = 45891
$ws.Cells.Item(77, 2).Value = 0.5

# Row 78: 23-Aug, no hours yet
$ws.Cells.Item(78, 1).Value = 45892

# Row 79: 24-Aug
$ws.Cells.Item(79, 1).Value = 45893

# Row 80: 25-Aug
$ws.Cells.Item(80, 1).Value = 45894

# Row 81: 26-Aug
$ws.Cells.Item(81, 1).Value = 45895

# Row 82: 27-Aug
$ws.Cells.Item(82, 1).Value = 45896

# Match date number format used by the rest of column A (d-mmm)
$ws.Range("A76:A82").NumberFormat = $ws.Range("A75").NumberFormat

# Update the selected cell to reflect where the user left off editing
$ws.Range("D78").Select()
